$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.668.98'
$ws.Range("E2").Value = '  +2.84%  '

$ws.Range("D3").Value = '4.030.36'
$ws.Range("E3").Value = '  +2.37%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''525.21'
$ws.Range("E5").Value = '  -1.19%  '

$ws.Range("D6").Value = '''148.40'
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("E7").Value = '  +0.56%  '

$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").Value = '''0.741'

$ws.Range("E10").Value = '  +1.27%  '

$ws.Range("E11").Value = '  -0.47%  '

$ws.Range("D12").Value = '''45.89'
$ws.Range("E12").Value = '  +6.73%  '

$ws.Range("D13").Value = '''10.81'
$ws.Range("E13").Value = '  +2.99%  '

$ws.Range("D14").Value = '4.666.02'
$ws.Range("E14").Value = '  +2.15%  '

$ws.Range("D15").Value = '4.015.64'
$ws.Range("E15").Value = '  +1.96%  '

$ws.Range("D16").Value = '''21.46'
$ws.Range("E16").Value = '  +7.73%  '

$ws.Range("D17").Value = '''14.29'
$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("D18").Value = '''1.22'
$ws.Range("E18").Value = '  +0.85%  '

$ws.Range("E19").Value = '  -1.66%  '

$ws.Range("D20").Value = '71.640.27'
$ws.Range("E20").Value = '  +2.94%  '

$ws.Range("D21").Value = '''443.17'
$ws.Range("E21").Value = '  +1.62%  '

$ws.Range("D22").Value = '''3.58'
$ws.Range("E22").Value = '  +4.65%  '

$ws.Range("D23").Value = '''95.11'
$ws.Range("E23").Value = '  +7.47%  '

$ws.Range("D24").Value = '''14.40'
$ws.Range("E24").Value = '  -1.43%  '

$ws.Range("D25").Value = '''12.33'
$ws.Range("E25").Value = '  +3.25%  '

$ws.Range("D26").Value = '''4.05'
$ws.Range("E26").Value = '  -0.58%  '

$ws.Range("E27").Value = '  +0.84%  '

$ws.Range("D28").Value = '''37.14'
$ws.Range("E28").Value = '  +0.97%  '

$ws.Range("D29").Value = '''13.64'
$ws.Range("E29").Value = '  +1.69%  '

$ws.Range("D30").Value = '''700.35'
$ws.Range("E30").Value = '  -0.74%  '

$ws.Range("E31").Value = '  +2.87%  '

$ws.Range("E32").Value = '  +1.54%  '

$ws.Range("E33").Value = '  +13.36%  '

$ws.Range("D34").Value = '''68.10'
$ws.Range("E34").Value = '  -0.60%  '

$ws.Range("D35").Value = '0.0₃0900'
$ws.Range("E35").Value = '  +3.21%  '

$ws.Range("E36").Value = '  +1.26%  '

$ws.Range("D37").Value = '''41.22'
$ws.Range("E37").Value = '  +1.53%  '

$ws.Range("D38").Value = '''0.158'
$ws.Range("E38").Value = '  +6.19%  '

$ws.Range("E39").Value = '  +17.06%  '

$ws.Range("E40").Value = '  +0.16%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  -0.10%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0492'
$ws.Range("E42").Value = '  +1.17%  '

$ws.Range("D43").Value = '''2.84'
$ws.Range("E43").Value = '  +0.30%  '

$ws.Range("D45").Value = '''3.52'
$ws.Range("E45").Value = '  +2.92%  '

$ws.Range("D46").Value = '''0.146'
$ws.Range("E46").Value = '  +2.55%  '

$ws.Range("E47").Value = '  -0.69%  '

$ws.Range("E48").Value = '  +6.82%  '

$ws.Range("B49").Value = 'LidoDAOToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D49").Value = '''3.38'
$ws.Range("E49").Value = '  +0.71%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '''0.000274'
$ws.Range("E50").Value = '  +15.04%  '

$ws.Range("D51").Value = '0.0₆0344'
$ws.Range("E51").Value = '  -5.58%  '
